# This is new modify code for datadriver framework.
#
# Reproduces the workbook edit:
#  - AF1 header relabelled "Status" -> "Result output"
#  - AF3 now carries the combined validation message
#    ("OrangeContact number is already in use. apple") and the stray AG3
#    duplicate of "Contact number is already in use." is removed
#  - A new data row (row 4) is appended, duplicating row 3's
#    BoatOwner record (PRASANNA / TARAI / SINGITALIA / ...), with its
#    own AF4 result message; the old row 4 (Gopal Mandal / Junusnagar)
#    and its stray AG4 cell are gone
#  - Page setup now specifies A4, portrait printing
#  - A few column widths grow to fit the new long text

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 (header) -------------------------------------------------
$ws.Range("AF1").Value = "Result output"

# --- Row 3: fold the two old messages into one on AF3, drop AG3 -----
$ws.Range("AF3").Value = "OrangeContact number is already in use. apple"
$ws.Range("AG3").ClearContents()

# --- Row 4: replace with a duplicate of row 3's record ---------------
# (same BoatOwner data as row 3, carrying row 3's number formats/styles
# along with it) and give it the same new AF result message.
$ws.Range("A3:AF3").Copy($ws.Range("A4:AF4"))
$ws.Range("AG4").ClearContents()

# --- Column widths: widen the columns that now hold longer text ------
$ws.Columns("L").ColumnWidth = 37.3
$ws.Columns("AD").ColumnWidth = 21
$ws.Columns("AF").ColumnWidth = 47.8

# --- Page setup: print on A4, portrait -------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Keep the active selection on the new result cell -----------------
$ws.Range("AF7").Select()

Write-Host "Applied datadriver framework update."
